# Apply the "Ran models for three Clavigralla populations" update:
#  - fill in the newly-computed r.model.h / r.model.f values for rows 2 and 4
#  - update the active selection to H3
#  - tighten the custom number format used by H4 from 0.000000 to 0.0000

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model-output values (G = r.model.h, H = r.model.f)
$ws.Range("G2").Value = 0.0153
$ws.Range("H2").Value = 0.0301

$ws.Range("G4").Value = 0.0182
$ws.Range("H4").Value = 0.0168

# H4 already carries a custom number format (six decimal places); narrow it
# down to four decimal places to match the refreshed precision.
$ws.Range("H4").NumberFormat = "0.0000"

# Move/restore the active cell selection to H3
$ws.Range("H3").Select()
